$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Twitchy Cat")

# --- Widen column D slightly (closest achievable to 18.3046875 in this engine) ---
$ws.Columns.Item(4).ColumnWidth = 17.5

# --- Core parameter edit: number of cats being built goes from 10 to 12 ---
$ws.Range("B2").Value = 12

# --- Cat Skeleton (row 5): price dropped, on-hand count corrected ---
$ws.Range("F5").Value = 17.670000000000002
$ws.Range("N5").Value = 1

# --- Hinge-ish item (row 6): price dropped ---
$ws.Range("F6").Value = 16.989999999999998

# --- Timer Relay (row 12): price increased ---
$ws.Range("F12").Value = 15.99

# --- row 9: orders + on-hand corrected ---
$ws.Range("K9").Value = 4
$ws.Range("N9").Value = 100

# --- Spring (row 15): price increased, on-hand corrected ---
$ws.Range("F15").Value = 10.09
$ws.Range("N15").Value = 7

# --- row 16: on-hand corrected ---
$ws.Range("N16").Value = 6

# --- row 18: orders corrected ---
$ws.Range("K18").Value = 6

# --- New grand-total cell for the Order Cost column ---
$ws.Range("Q26").Formula = "=SUM(Q5:Q24)"
$ws.Range("Q26").NumberFormat = $ws.Range("Q5").NumberFormat

# --- Turn the remaining bare source-URL cells into real hyperlinks ---
function Add-UrlHyperlink($cellRef) {
    $cell = $ws.Range($cellRef)
    $url = $cell.Value
    $ws.Hyperlinks.Add($cell, $url)
    $cell.Style = "Hyperlink"
}

Add-UrlHyperlink "D5"
Add-UrlHyperlink "D6"
Add-UrlHyperlink "D7"
Add-UrlHyperlink "D13"
Add-UrlHyperlink "D14"
Add-UrlHyperlink "D15"
Add-UrlHyperlink "D16"

# --- Restore the selection to where the author left off ---
$ws.Range("D18").Select()
